$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 25.99000000000062
$ws.Range("H2").Value = 0.03842527030919274
$ws.Range("I2").Value = 0.03842527030919274
$ws.Range("L2").Value = 7.141271080794628
$ws.Range("M2").Value = "[-0.07925369504240187, 14.361795856631657]"
$ws.Range("N2").Value = 0.05246022008588613
$ws.Range("O2").Value = 0.05246022008588613
$ws.Range("P2").Value = -0.8176317216550011
$ws.Range("Q2").Value = "[-3.1195794918529276, 1.4843160485429254]"
$ws.Range("R2").Value = 0.4780625111959962
$ws.Range("S2").Value = 0.4780625111959962
$ws.Range("T2").Value = 12.35775680271132
$ws.Range("U2").Value = "[8.405081507977929, 16.310432097444707]"
$ws.Range("V2").Value = [double]"1.130209545951999e-07"
$ws.Range("W2").Value = [double]"1.130209545951999e-07"
$ws.Range("X2").Value = 3.382082082082164
$ws.Range("Y2").Value = -6.139779779779927
$ws.Range("Z2").Value = 12.90394394394425

# Row 3
$ws.Range("F3").Value = 25.99000000000062
$ws.Range("H3").Value = 0.2064319267337758
$ws.Range("I3").Value = 0.2064319267337758
$ws.Range("L3").Value = 6.51416860119548
$ws.Range("M3").Value = "[-3.0829289790402674, 16.111266181431226]"
$ws.Range("N3").Value = 0.1783841777821056
$ws.Range("O3").Value = 0.1783841777821056
$ws.Range("P3").Value = -2.012631930227695
$ws.Range("Q3").Value = "[-4.956106128185699, 0.9308422677303092]"
$ws.Range("R3").Value = 0.175273598113558
$ws.Range("S3").Value = 0.175273598113558
$ws.Range("T3").Value = 16.49033631332604
$ws.Range("U3").Value = "[11.322212282001196, 21.658460344650884]"
$ws.Range("V3").Value = [double]"7.250254374646659e-08"
$ws.Range("W3").Value = [double]"7.250254374646659e-08"
$ws.Range("X3").Value = 8.325125125125322
$ws.Range("Y3").Value = -3.850370370370465
$ws.Range("Z3").Value = 20.50062062062111

# Row 4
$ws.Range("F4").Value = 25.99000000000062
$ws.Range("H4").Value = 0.3011618927652202
$ws.Range("I4").Value = 0.3011618927652202
$ws.Range("L4").Value = 4.185561313134684
$ws.Range("M4").Value = "[-2.6761543838370656, 11.047277010106434]"
$ws.Range("N4").Value = 0.2256166798333348
$ws.Range("O4").Value = 0.2256166798333348
$ws.Range("P4").Value = 1.37739497724958
$ws.Range("Q4").Value = "[-1.7296055650394244, 4.484395519538584]"
$ws.Range("R4").Value = 0.3766652765878722
$ws.Range("S4").Value = 0.3766652765878722
$ws.Range("T4").Value = 14.03789296946483
$ws.Range("U4").Value = "[10.136958884625038, 17.938827054304618]"
$ws.Range("V4").Value = [double]"4.376023099439408e-09"
$ws.Range("W4").Value = [double]"4.376023099439408e-09"
$ws.Range("X4").Value = 20.29249249249298
$ws.Range("Y4").Value = 7.440580580580757
$ws.Range("Z4").Value = 33.14440440440519

# Row 5
$ws.Range("B5").Value = 1
$ws.Range("F5").Value = 25.99000000000062
$ws.Range("H5").Value = 0.01459907808970407
$ws.Range("I5").Value = 0.01459907808970407
$ws.Range("L5").Value = 9.458929849244036
$ws.Range("M5").Value = "[1.012901441902546, 17.904958256585527]"
$ws.Range("N5").Value = 0.02900149984792888
$ws.Range("O5").Value = 0.02900149984792888
$ws.Range("P5").Value = -2.339684618889696
$ws.Range("Q5").Value = "[-3.069263693597235, -1.6101055441821561]"
$ws.Range("R5").Value = [double]"6.487163228818815e-08"
$ws.Range("S5").Value = [double]"6.487163228818815e-08"
$ws.Range("T5").Value = 12.56652179813175
$ws.Range("U5").Value = "[8.238899080893091, 16.8941445153704]"
$ws.Range("V5").Value = [double]"5.237367997068532e-07"
$ws.Range("W5").Value = [double]"5.237367997068532e-07"
$ws.Range("X5").Value = 9.677957957958188
$ws.Range("Y5").Value = 6.660100100100259
$ws.Range("Z5").Value = 12.69581581581612

# Row 6
$ws.Range("F6").Value = 25.99000000000062
$ws.Range("H6").Value = 0.3771003798725021
$ws.Range("I6").Value = 0.3771003798725021
$ws.Range("L6").Value = 5.364018756012653
$ws.Range("M6").Value = "[-4.098815646394037, 14.826853158419343]"
$ws.Range("N6").Value = 0.2596190439491675
$ws.Range("O6").Value = 0.2596190439491675
$ws.Range("P6").Value = -3.056684744033312
$ws.Range("Q6").Value = "[-6.169974761104277, 0.05660527303765406]"
$ws.Range("R6").Value = 0.05413159631972397
$ws.Range("S6").Value = 0.05413159631972397
$ws.Range("T6").Value = 13.4239734889609
$ws.Range("U6").Value = "[8.215880327175253, 18.63206665074655]"
$ws.Range("V6").Value = [double]"4.851388203608309e-06"
$ws.Range("W6").Value = [double]"4.851388203608309e-06"
$ws.Range("X6").Value = 12.64378378378409
$ws.Range("Y6").Value = -0.2341441441441479
$ws.Range("Z6").Value = 25.52171171171232

# Row 7
$ws.Range("B7").Value = 0
$ws.Range("F7").Value = 25.99000000000062
$ws.Range("H7").Value = 0.4295349290348131
$ws.Range("I7").Value = 0.4295349290348131
$ws.Range("L7").Value = 4.375674926031047
$ws.Range("M7").Value = "[-4.695237835377536, 13.446587687439628]"
$ws.Range("N7").Value = 0.3364562214721061
$ws.Range("O7").Value = 0.3364562214721061
$ws.Range("P7").Value = 2.685605731897581
$ws.Range("Q7").Value = "[-0.4339737599553466, 5.805185223750509]"
$ws.Range("R7").Value = 0.08978010843072726
$ws.Range("S7").Value = 0.08978010843072726
$ws.Range("T7").Value = 13.3440619336396
$ws.Range("U7").Value = "[8.58735822443446, 18.10076564284475]"
$ws.Range("V7").Value = [double]"1.029161593946171e-06"
$ws.Range("W7").Value = [double]"1.029161593946171e-06"
$ws.Range("X7").Value = 14.88116116116152
$ws.Range("Y7").Value = 1.977217217217264
$ws.Range("Z7").Value = 27.78510510510577

# Row 8
$ws.Range("F8").Value = 25.99000000000062
$ws.Range("H8").Value = 0.4578176744023924
$ws.Range("I8").Value = 0.4578176744023924
$ws.Range("L8").Value = 4.369074292875869
$ws.Range("M8").Value = "[-5.298696830459594, 14.036845416211332]"
$ws.Range("N8").Value = 0.3675580774203018
$ws.Range("O8").Value = 0.3675580774203018
$ws.Range("P8").Value = 2.723342580589351
$ws.Range("Q8").Value = "[-0.40881586082749877, 5.855501022006201]"
$ws.Range("R8").Value = 0.08672242495145355
$ws.Range("S8").Value = 0.08672242495145355
$ws.Range("T8").Value = 13.45925770166922
$ws.Range("U8").Value = "[8.534808535251624, 18.383706868086808]"
$ws.Range("V8").Value = [double]"1.685508316251472e-06"
$ws.Range("W8").Value = [double]"1.685508316251472e-06"
$ws.Range("X8").Value = 14.72506506506542
$ws.Range("Y8").Value = 1.769089089089134
$ws.Range("Z8").Value = 27.68104104104169

# Row 9
$ws.Range("F9").Value = 23.9000000000003
$ws.Range("H9").Value = 0.1617660496307021
$ws.Range("I9").Value = 0.1617660496307021
$ws.Range("L9").Value = 6.111705916867909
$ws.Range("M9").Value = "[-1.991519845642105, 14.214931679377923]"
$ws.Range("N9").Value = 0.1357327523632554
$ws.Range("O9").Value = 0.1357327523632554
$ws.Range("P9").Value = 1.427710775505271
$ws.Range("Q9").Value = "[-0.861658045128733, 3.717079596139275]"
$ws.Range("R9").Value = 0.2155807568725971
$ws.Range("S9").Value = 0.2155807568725971
$ws.Range("T9").Value = 15.01574619787992
$ws.Range("U9").Value = "[10.560997279713861, 19.47049511604598]"
$ws.Range("V9").Value = [double]"2.095956186387582e-08"
$ws.Range("W9").Value = [double]"2.095956186387582e-08"
$ws.Range("X9").Value = 18.4692692692695
$ws.Range("Y9").Value = 9.76096096096108
$ws.Range("Z9").Value = 27.17757757757792

# Row 10
$ws.Range("F10").Value = 23.9000000000003
$ws.Range("H10").Value = 0.1423436750893533
$ws.Range("I10").Value = 0.1423436750893533
$ws.Range("L10").Value = 5.830963001689888
$ws.Range("M10").Value = "[-2.189852889549062, 13.851778892928838]"
$ws.Range("N10").Value = 0.1500878745694232
$ws.Range("O10").Value = 0.1500878745694232
$ws.Range("P10").Value = 2.283079345852042
$ws.Range("Q10").Value = "[-0.7987632973091161, 5.3649219890132]"
$ws.Range("R10").Value = 0.1426559752041294
$ws.Range("S10").Value = 0.1426559752041294
$ws.Range("T10").Value = 11.66064321722349
$ws.Range("U10").Value = "[7.548348878714288, 15.772937555732687]"
$ws.Range("V10").Value = [double]"8.36583017971293e-07"
$ws.Range("W10").Value = [double]"8.36583017971293e-07"
$ws.Range("X10").Value = 15.21561561561581
$ws.Range("Y10").Value = 3.492892892892938
$ws.Range("Z10").Value = 26.93833833833867

# Row 11
$ws.Range("F11").Value = 23.9000000000003
$ws.Range("H11").Value = 0.06918484791380508
$ws.Range("I11").Value = 0.06918484791380508
$ws.Range("L11").Value = 7.720715015783436
$ws.Range("M11").Value = "[-0.9756989038542203, 16.417128935421093]"
$ws.Range("N11").Value = 0.08048937885531782
$ws.Range("O11").Value = 0.08048937885531782
$ws.Range("P11").Value = 2.459184639746965
$ws.Range("Q11").Value = "[-0.3018947895341544, 5.220264069028084]"
$ws.Range("R11").Value = 0.07955268140690697
$ws.Range("S11").Value = 0.07955268140690697
$ws.Range("T11").Value = 14.690676985093
$ws.Range("U11").Value = "[10.07331660579623, 19.308037364389776]"
$ws.Range("V11").Value = [double]"7.723257167135955e-08"
$ws.Range("W11").Value = [double]"7.723257167135955e-08"
$ws.Range("X11").Value = 14.54574574574593
$ws.Range("Y11").Value = 4.043143143143196
$ws.Range("Z11").Value = 25.04834834834866

# Row 12
$ws.Range("F12").Value = 23.9000000000003
$ws.Range("H12").Value = 0.2467553678210851
$ws.Range("I12").Value = 0.2467553678210851
$ws.Range("L12").Value = 6.331573838542631
$ws.Range("M12").Value = "[-3.487706005034365, 16.150853682119628]"
$ws.Range("N12").Value = 0.2006576783646885
$ws.Range("O12").Value = 0.2006576783646885
$ws.Range("P12").Value = 2.673026782333658
$ws.Range("Q12").Value = "[-0.44655270951926873, 5.792606274186586]"
$ws.Range("R12").Value = 0.09124593069532816
$ws.Range("S12").Value = 0.09124593069532816
$ws.Range("T12").Value = 13.10626876755355
$ws.Range("U12").Value = "[7.795430865578904, 18.4171066695282]"
$ws.Range("V12").Value = [double]"1.014881672944945e-05"
$ws.Range("W12").Value = [double]"1.014881672944945e-05"
$ws.Range("X12").Value = 13.7323323323325
$ws.Range("Y12").Value = 1.866066066066089
$ws.Range("Z12").Value = 25.59859859859891

# Row 13
$ws.Range("F13").Value = 23.9000000000003
$ws.Range("H13").Value = 0.1720133354592219
$ws.Range("I13").Value = 0.1720133354592219
$ws.Range("L13").Value = 6.200054594464074
$ws.Range("M13").Value = "[-2.541621438356385, 14.941730627284533]"
$ws.Range("N13").Value = 0.1600503822395907
$ws.Range("O13").Value = 0.1600503822395907
$ws.Range("P13").Value = 2.333395144107734
$ws.Range("Q13").Value = "[-0.798763297309117, 5.465553585524585]"
$ws.Range("R13").Value = 0.1404778346629785
$ws.Range("S13").Value = 0.1404778346629785
$ws.Range("T13").Value = 15.24601560595342
$ws.Range("U13").Value = "[10.638187017146791, 19.85384419476005]"
$ws.Range("V13").Value = [double]"3.213746868624412e-08"
$ws.Range("W13").Value = [double]"3.213746868624412e-08"
$ws.Range("X13").Value = 15.02422422422441
$ws.Range("Y13").Value = 3.110110110110147
$ws.Range("Z13").Value = 26.93833833833867

# Row 14
$ws.Range("F14").Value = 23.9000000000003
$ws.Range("H14").Value = 0.5108404923687053
$ws.Range("I14").Value = 0.5108404923687053
$ws.Range("L14").Value = 3.659067213076744
$ws.Range("M14").Value = "[-4.054185812435488, 11.372320238588976]"
$ws.Range("N14").Value = 0.3444477344878374
$ws.Range("O14").Value = 0.3444477344878374
$ws.Range("P14").Value = 1.578658170272348
$ws.Range("Q14").Value = "[-1.553500271144503, 4.710816611689199]"
$ws.Range("R14").Value = 0.3154660619721124
$ws.Range("S14").Value = 0.3154660619721124
$ws.Range("T14").Value = 11.41607621422065
$ws.Range("U14").Value = "[6.9654279629558715, 15.86672446548543]"
$ws.Range("V14").Value = [double]"5.278300961419902e-06"
$ws.Range("W14").Value = [double]"5.278300961419902e-06"
$ws.Range("X14").Value = 17.89509509509532
$ws.Range("Y14").Value = 5.980980980981053
$ws.Range("Z14").Value = 29.80920920920958
